$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2718
$ws.Range("F5").Value = 218
$ws.Range("F6").Value = 523
$ws.Range("F7").Value = 1301
$ws.Range("F8").Value = 615
$ws.Range("F9").Value = 331
$ws.Range("F12").Value = 397
$ws.Range("F13").Value = 6008
$ws.Range("F14").Value = 106
$ws.Range("F16").Value = 1864
$ws.Range("F17").Value = 4478
$ws.Range("F21").Value = 5219
$ws.Range("F22").Value = 6783
$ws.Range("F24").Value = 1073
$ws.Range("F25").Value = 730
$ws.Range("F26").Value = 3916
$ws.Range("F27").Value = 530
$ws.Range("F29").Value = 211
$ws.Range("F31").Value = 1028
$ws.Range("F32").Value = 1468
$ws.Range("F33").Value = 527
$ws.Range("F34").Value = 636
$ws.Range("F35").Value = 1652
$ws.Range("F36").Value = 225
$ws.Range("F37").Value = 1828
$ws.Range("C38").Value = "杭州·梦漫星河动漫展"
$ws.Range("D38").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E38").Value = "2024.08.03 10:00-08.04 17:00"
$ws.Range("F38").Value = 1202
$ws.Range("G38").Value = 68
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=82836"
$ws.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202403/VFfQUJdD1711700169290.jpeg"
$ws.Range("C39").Value = "杭州·梦漫星河动漫嘉年华·赵路专场"
$ws.Range("E39").Value = "2024.08.04 11:40-08.04 17:00"
$ws.Range("F39").Value = "已售罄"
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=86221"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202405/2padflbr1716372780297.jpeg"
$ws.Range("C40").Value = "杭州·原神X星铁X绝区零only"
$ws.Range("D40").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E40").Value = "2024.08.10 10:00-08.10 17:00"
$ws.Range("F40").Value = 664
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
$ws.Range("C41").Value = "杭州·造梦探险家城堡二次元同好会"
$ws.Range("D41").Value = "大岭山路156号 爱丽芬城堡"
$ws.Range("E41").Value = "2024.08.10 10:00-08.10 22:00"
$ws.Range("F41").Value = 111
$ws.Range("G41").Value = 38
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=86432"
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202405/xWUy30Ns1716783723057.jpeg"
$ws.Range("C42").Value = "【会员购严选】杭州·首届次元之门动漫游戏博览会"
$ws.Range("D42").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E42").Value = "2024.08.17 10:00-08.18 17:30"
$ws.Range("F42").Value = 265
$ws.Range("G42").Value = 75
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=87065"
$ws.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202406/wrxORgrP1717593610187.jpeg"
$ws.Range("F43").Value = 3581
$ws.Range("F44").Value = 145
$ws.Range("F46").Value = 431
$ws.Range("F47").Value = 15
$ws.Range("F48").Value = 73
$ws.Range("F49").Value = 3925
foreach ($dc in @(@("B39","2024-08-04"), @("B40","2024-08-10"), @("B42","2024-08-17"))) {
    $ws.Range($dc[0]).NumberFormat = "@"
    $ws.Range($dc[0]).Value = $dc[1]
    $ws.Range($dc[0]).Style = "Normal"
}

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 1246
$ws.Range("F15").Value = 19
$ws.Range("F20").Value = 1
$ws.Range("F28").Value = 23

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4233

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 4233
$ws.Range("F3").Value = 2718
$ws.Range("F6").Value = 1246
$ws.Range("F8").Value = 218
$ws.Range("F9").Value = 523
$ws.Range("F11").Value = 1301
$ws.Range("F13").Value = 615
$ws.Range("F14").Value = 331
$ws.Range("F17").Value = 397
$ws.Range("F20").Value = 1864
$ws.Range("F21").Value = 4478
$ws.Range("F22").Value = 5219
$ws.Range("F23").Value = 5219
$ws.Range("F25").Value = 1073
$ws.Range("F26").Value = 730
$ws.Range("F27").Value = 3916
$ws.Range("F28").Value = 530
$ws.Range("F31").Value = 1028
$ws.Range("F32").Value = 1468
$ws.Range("F33").Value = 527
$ws.Range("F34").Value = 636
$ws.Range("F35").Value = 1652
$ws.Range("F36").Value = 1828
$ws.Range("F37").Value = 1
$ws.Range("F39").Value = 664
$ws.Range("F41").Value = 111
$ws.Range("F43").Value = 3581
$ws.Range("F44").Value = 23
$ws.Range("F45").Value = 145
$ws.Range("F47").Value = 431
$ws.Range("F48").Value = 74
$ws.Range("F50").Value = 3925
